$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Fix normals problem on cylinder" task (row 2) has been fixed, so
# remove its row entirely - this shifts every following row up by one
# and drops its now-unused shared string / style.
$ws.Rows.Item(2).Delete()

# The "Vertex welding in model compiler" task (now row 3, after the
# deletion above) needs to be reprioritized to the bottom of the list.
# Remember its values, shift the rows below it up, then write the
# remembered values into the now-freed last row.
$taskValue = $ws.Range("A3").Value2
$estimateValue = $ws.Range("B3").Value2

for ($row = 3; $row -le 5; $row++) {
    $nextRow = $row + 1
    $ws.Range("A$row").Value2 = $ws.Range("A$nextRow").Value2
    $ws.Range("B$row").Value2 = $ws.Range("B$nextRow").Value2
}

$ws.Range("A6").Value2 = $taskValue
$ws.Range("B6").Value2 = $estimateValue

# Leave the same selection state as after performing the row delete by hand.
$ws.Rows.Item(2).Select() | Out-Null
